$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 421
$ws1.Range("F5").Value = 1252
$ws1.Range("F7").Value = 7523
$ws1.Range("F11").Value = 8189
$ws1.Range("F14").Value = 5589
$ws1.Range("F16").Value = 2551
$ws1.Range("F18").Value = 4576
$ws1.Range("F19").Value = 322
$ws1.Range("F21").Value = 89
$ws1.Range("F24").Value = 2146
$ws1.Range("F26").Value = 2737
$ws1.Range("F28").Value = 313
$ws1.Range("F29").Value = 110
$ws1.Range("F30").Value = 261
$ws1.Range("F34").Value = 1610
$ws1.Range("F37").Value = 2566

# Sheet "演出" (sheet2)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 104
$ws2.Range("F4").Value = 35
$ws2.Range("F6").Value = 32
$ws2.Range("F8").Value = 100

# Sheet "本地生活" (sheet3)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 262
$ws3.Range("F3").Value = 1298

# Sheet "全部类型" (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 262
$ws4.Range("F4").Value = 1298
$ws4.Range("F5").Value = 421
$ws4.Range("F6").Value = 1252
$ws4.Range("F7").Value = 461
$ws4.Range("F8").Value = 7523
$ws4.Range("F12").Value = 8189
$ws4.Range("F15").Value = 5589
$ws4.Range("F17").Value = 2551
$ws4.Range("F19").Value = 4576
$ws4.Range("F20").Value = 322
$ws4.Range("F22").Value = 89
$ws4.Range("F25").Value = 104
$ws4.Range("F27").Value = 2146
$ws4.Range("F29").Value = 2737
$ws4.Range("F31").Value = 313
$ws4.Range("F32").Value = 110
$ws4.Range("F33").Value = 261
$ws4.Range("F34").Value = 35
$ws4.Range("F39").Value = 32
$ws4.Range("F40").Value = 1610
$ws4.Range("F43").Value = 2566
$ws4.Range("F49").Value = 100

$wb.Save()
